# The deck's theme (ppt/theme/theme1.xml, "Integral") is swapped for the
# stock Office theme's colour scheme (the colours that lived in
# ppt/theme/theme2.xml, "Office Theme", used by the notes master) while the
# font scheme / format scheme (identical between the two theme parts already)
# stay as-is.
#
# MsoThemeColorSchemeIndex order exposed via ThemeColorScheme:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB() values below are R + G*256 + B*65536 for the target "Office" palette:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6 accent1=5B9BD5 accent2=ED7D31
#   accent3=A5A5A5 accent4=FFC000 accent5=4472C4 accent6=70AD47
#   hlink=0563C1 folHlink=954F72

$p = $ppt.ActivePresentation

$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
